# Update cryptos list with latest prices / volume changes (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading ' forces Excel to store the literal text (keeps '1.00', '0.0620', etc.
# instead of silently coercing to a Double) without touching NumberFormat.

$ws.Range("D2").Value = "'36.622.69"
$ws.Range("E2").Value = "'  +3.84%  "

$ws.Range("D3").Value = "'1.913.14"
$ws.Range("E3").Value = "'  +1.89%  "

$ws.Range("E4").Value = "'  -0.11%  "

$ws.Range("B5").Value = "'XRP"
$ws.Range("C5").Value = "'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D5").Value = "'0.700"
$ws.Range("E5").Value = "'  +3.35%  "

$ws.Range("B6").Value = "'BNB"
$ws.Range("C6").Value = "'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D6").Value = "'248.82"
$ws.Range("E6").Value = "'  +1.70%  "

$ws.Range("E7").Value = "'  -0.04%  "

$ws.Range("D8").Value = "'44.54"
$ws.Range("E8").Value = "'  +1.76%  "

$ws.Range("E9").Value = "'  +3.32%  "

$ws.Range("D10").Value = "'58.09"
$ws.Range("E10").Value = "'  +8.89%  "

$ws.Range("D11").Value = "'0.0763"
$ws.Range("E11").Value = "'  +3.30%  "

$ws.Range("E12").Value = "'  +2.48%  "

$ws.Range("D13").Value = "'14.49"
$ws.Range("E13").Value = "'  +7.07%  "

$ws.Range("E14").Value = "'  +5.66%  "

$ws.Range("D15").Value = "'2.191.78"
$ws.Range("E15").Value = "'  +1.87%  "

$ws.Range("D16").Value = "'5.11"
$ws.Range("E16").Value = "'  +4.02%  "

$ws.Range("D17").Value = "'1.916.01"
$ws.Range("E17").Value = "'  +2.06%  "

$ws.Range("D18").Value = "'36.669.00"
$ws.Range("E18").Value = "'  +3.94%  "

$ws.Range("D19").Value = "'74.42"
$ws.Range("E19").Value = "'  +1.53%  "

$ws.Range("D20").Value = "'0.0₃0862"
$ws.Range("E20").Value = "'  +5.15%  "

$ws.Range("D21").Value = "'250.08"
$ws.Range("E21").Value = "'  +2.48%  "

$ws.Range("E22").Value = "'  +4.21%  "

$ws.Range("E23").Value = "'  +3.39%  "

$ws.Range("E24").Value = "'  -2.76%  "

$ws.Range("D25").Value = "'0.999"
$ws.Range("E25").Value = "'  -0.07%  "

$ws.Range("D26").Value = "'2.20"
$ws.Range("E26").Value = "'  +2.43%  "

$ws.Range("D27").Value = "'168.59"
$ws.Range("E27").Value = "'  +2.40%  "

$ws.Range("E28").Value = "'  +2.58%  "

$ws.Range("E29").Value = "'  +2.45%  "

$ws.Range("E30").Value = "'  +1.70%  "

$ws.Range("E31").Value = "'  +6.82%  "

$ws.Range("D32").Value = "'0.0620"
$ws.Range("E32").Value = "'  +5.03%  "

$ws.Range("D33").Value = "'4.34"
$ws.Range("E33").Value = "'  +4.43%  "

$ws.Range("B34").Value = "'Kaspa"
$ws.Range("C34").Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D34").Value = "'0.0901"
$ws.Range("E34").Value = "'  +24.12%  "

$ws.Range("B35").Value = "'WEMIXToken"
$ws.Range("C35").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D35").Value = "'1.92"
$ws.Range("E35").Value = "'  +5.69%  "

$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "'  -0.07%  "

$ws.Range("E37").Value = "'  +7.12%  "

$ws.Range("D38").Value = "'0.880"
$ws.Range("E38").Value = "'  +3.67%  "

$ws.Range("D39").Value = "'17.83"
$ws.Range("E39").Value = "'  +52.38%  "

$ws.Range("E40").Value = "'  +4.88%  "

$ws.Range("D41").Value = "'106.29"
$ws.Range("E41").Value = "'  +10.41%  "

$ws.Range("D42").Value = "'0.0227"
$ws.Range("E42").Value = "'  +4.37%  "

$ws.Range("D43").Value = "'17.47"
$ws.Range("E43").Value = "'  +0.97%  "

$ws.Range("E44").Value = "'  +23.46%  "

$ws.Range("D45").Value = "'1.11"
$ws.Range("E45").Value = "'  +3.37%  "

$ws.Range("D46").Value = "'1.343.69"
$ws.Range("E46").Value = "'  +3.01%  "

$ws.Range("E47").Value = "'  +0.07%  "

$ws.Range("D48").Value = "'0.0813"
$ws.Range("E48").Value = "'  +1.87%  "

$ws.Range("E49").Value = "'  +2.02%  "

$ws.Range("D50").Value = "'6.39"
$ws.Range("E50").Value = "'  +2.00%  "

$ws.Range("D51").Value = "'43.46"
$ws.Range("E51").Value = "'  +3.65%  "
